$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.246.94"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.214.88"
$ws.Range("E3").Value = "  +2.60%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.73"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.39"
$ws.Range("E6").Value = "  +4.77%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.531"
$ws.Range("E8").Value = "  +3.33%  "
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("E10").Value = "  +4.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.435"
$ws.Range("E11").Value = "  +3.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.771.45"
$ws.Range("E12").Value = "  +2.61%  "
$ws.Range("E13").Value = "  -1.07%  "
$ws.Range("E14").Value = "  +3.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.15"
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "60.305.55"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.208.35"
$ws.Range("E17").Value = "  +1.88%  "
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.20"
$ws.Range("E19").Value = "  +1.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.37"
$ws.Range("E20").Value = "  +2.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "383.24"
$ws.Range("E21").Value = "  +2.22%  "
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("E23").Value = "  +3.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.26"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("E25").Value = "  +11.81%  "
$ws.Range("E26").Value = "  +2.30%  "
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0914"
$ws.Range("E28").Value = "  +4.18%  "
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.51"
$ws.Range("E30").Value = "  +3.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.47"
$ws.Range("E31").Value = "  +5.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.21"
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("E33").Value = "  +5.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.63"
$ws.Range("E34").Value = "  +6.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "157.00"
$ws.Range("E35").Value = "  -2.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.36"
$ws.Range("E36").Value = "  +0.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.787.35"
$ws.Range("E37").Value = "  +5.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.88"
$ws.Range("E38").Value = "  +1.54%  "
$ws.Range("E39").Value = "  +4.80%  "
$ws.Range("E40").Value = "  +0.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.28"
$ws.Range("E41").Value = "  +0.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.94"
$ws.Range("E42").Value = "  +4.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.726"
$ws.Range("E43").Value = "  +3.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0286"
$ws.Range("E44").Value = "  +4.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.259.21"
$ws.Range("E45").Value = "  +2.55%  "
$ws.Range("E46").Value = "  +3.55%  "
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.19"
$ws.Range("E48").Value = "  -0.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.808"
$ws.Range("E49").Value = "  +8.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.80"
$ws.Range("E50").Value = "  +3.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  +0.06%  "
